# Edit script: Log Week 16 stats and perform season sim from Week 17
# Panthers Players Data - Rushing & Receiving sheets

$wb = $excel.ActiveWorkbook

# ---------- Rushing sheet ----------
$wsR = $wb.Worksheets.Item("Rushing")

# Current sheet has 10 data rows (rows 2-11). Target has 12 data rows (rows 2-13).
# Copy the column-A index-cell formatting onto the two new trailing rows before writing values,
# so the new cells pick up the same bold/bordered style as the existing index column.
$wsR.Cells.Item(11,1).Copy()
$wsR.Cells.Item(12,1).PasteSpecial(-4122)
$wsR.Cells.Item(11,1).Copy()
$wsR.Cells.Item(13,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rushingRows = @(
    ,@(0, "C.Newton", 24, 23, 17, 10)
    ,@(1, "P.Walker", 1, 3, 2, 0)
    ,@(2, "S.Darnold", 1, 0, 0, 0)
    ,@(3, "C.Hubbard", 35, 17, 7, 9)
    ,@(4, "R.Freeman", 6, 3, 1, 0)
    ,@(5, "A.Abdullah", 14, 13, 4, 4)
    ,@(6, "R.Bonnafon", 1, 0, 0, 0)
    ,@(7, "R.Anderson", 2, 0, 0, 1)
    ,@(8, "Dj.Moore", 1, 2, 1, 0)
    ,@(9, "S.Smith", 0, 1, 0, 0)
    ,@(10, "T.Tremble", 0, 1, 0, 1)
    ,@(11, "I.Thomas", 1, 0, 0, 0)
)

foreach ($row in $rushingRows) {
    $r = [int]$row[0] + 2
    $wsR.Cells.Item($r,1).Value = $row[0]
    $wsR.Cells.Item($r,2).Value = $row[1]
    $wsR.Cells.Item($r,3).Value = $row[2]
    $wsR.Cells.Item($r,4).Value = $row[3]
    $wsR.Cells.Item($r,5).Value = $row[4]
    $wsR.Cells.Item($r,6).Value = $row[5]
}

# ---------- Receiving sheet ----------
$wsC = $wb.Worksheets.Item("Receiving")

# Current sheet has 14 data rows (rows 2-15). Target has 15 data rows (rows 2-16).
$wsC.Cells.Item(15,1).Copy()
$wsC.Cells.Item(16,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$receivingRows = @(
    ,@(0, "C.Hubbard", 17, 10, 2, 2, 0, 0)
    ,@(1, "R.Freeman", 6, 3, 0, 0, 0, 0)
    ,@(2, "A.Abdullah", 31, 21, 3, 1, 4, 2)
    ,@(3, "R.Bonnafon", 2, 2, 0, 0, 1, 1)
    ,@(4, "R.Anderson", 77, 40, 21, 3, 6, 2)
    ,@(5, "Dj.Moore", 107, 69, 38, 14, 13, 6)
    ,@(6, "T.Marshall", 22, 14, 6, 1, 3, 2)
    ,@(7, "B.Zylstra", 17, 14, 6, 4, 3, 2)
    ,@(8, "S.Smith", 8, 5, 2, 1, 1, 1)
    ,@(9, "K.Kirkwood", 4, 3, 2, 0, 1, 0)
    ,@(10, "W.Snead", 2, 1, 1, 0, 0, 0)
    ,@(11, "A.Erickson", 2, 1, 0, 0, 0, 0)
    ,@(12, "T.Tremble", 24, 17, 6, 1, 3, 3)
    ,@(13, "I.Thomas", 24, 13, 3, 2, 4, 2)
    ,@(14, "C.Thompson", 1, 0, 0, 0, 0, 0)
)

foreach ($row in $receivingRows) {
    $r = [int]$row[0] + 2
    $wsC.Cells.Item($r,1).Value = $row[0]
    $wsC.Cells.Item($r,2).Value = $row[1]
    $wsC.Cells.Item($r,3).Value = $row[2]
    $wsC.Cells.Item($r,4).Value = $row[3]
    $wsC.Cells.Item($r,5).Value = $row[4]
    $wsC.Cells.Item($r,6).Value = $row[5]
    $wsC.Cells.Item($r,7).Value = $row[6]
    $wsC.Cells.Item($r,8).Value = $row[7]
}

